$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D, E
$ws.Range("D2").Value = "42.099.28"
$ws.Range("E2").Value = "  -2.07%  "

# Row 3: update D, E
$ws.Range("D3").Value = "2.265.18"
$ws.Range("E3").Value = "  -3.11%  "

# Row 4: update E
$ws.Range("E4").Value = "  +0.01%  "

# Row 5: update D
$ws.Range("D5").Value = "298.07"

# Row 6: update D, E
$ws.Range("D6").Value = "94.04"
$ws.Range("E6").Value = "  -7.38%  "

# Row 7: update D, E
$ws.Range("D7").Value = "0.497"
$ws.Range("E7").Value = "  -2.86%  "

# Row 8: update E
$ws.Range("E8").Value = "  +0.08%  "

# Row 9: update E
$ws.Range("E9").Value = "  -4.24%  "

# Row 10: update D, E
$ws.Range("D10").Value = "32.97"
$ws.Range("E10").Value = "  -5.70%  "

# Row 11: update E
$ws.Range("E11").Value = "  -1.34%  "

# Row 12: update D, E
$ws.Range("D12").Value = "48.33"
$ws.Range("E12").Value = "  -7.49%  "

# Row 13: update E
$ws.Range("E13").Value = "  -0.06%  "

# Row 14: update D, E
$ws.Range("D14").Value = "6.65"
$ws.Range("E14").Value = "  -2.77%  "

# Row 15: update B, C, D, E
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.616.00"
$ws.Range("E15").Value = "  -3.20%  "

# Row 16: update B, C, D, E
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "15.51"
$ws.Range("E16").Value = "  -2.49%  "

# Row 17: update D, E
$ws.Range("D17").Value = "2.261.93"
$ws.Range("E17").Value = "  -2.53%  "

# Row 18: update E
$ws.Range("E18").Value = "  -5.80%  "

# Row 19: update D, E
$ws.Range("D19").Value = "42.089.06"
$ws.Range("E19").Value = "  -1.92%  "

# Row 20: update D, E
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  -2.44%  "

# Row 21: update E
$ws.Range("E21").Value = "  -3.50%  "

# Row 22: update E
$ws.Range("E22").Value = "  -3.57%  "

# Row 23: update D, E
$ws.Range("D23").Value = "66.67"
$ws.Range("E23").Value = "  -1.92%  "

# Row 24: update D, E
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.10"
$ws.Range("E24").Value = "  -1.78%  "

# Row 25: update E
$ws.Range("E25").Value = "  -4.73%  "

# Row 26: update D, E
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.03%  "

# Row 27: update E
$ws.Range("E27").Value = "  -4.40%  "

# Row 28: update D, E
$ws.Range("D28").Value = "23.78"
$ws.Range("E28").Value = "  -6.60%  "

# Row 29: update E
$ws.Range("E29").Value = "  -1.11%  "

# Row 30: update D, E
$ws.Range("D30").Value = "167.56"
$ws.Range("E30").Value = "  +3.99%  "

# Row 31: update D, E
$ws.Range("D31").Value = "33.61"
$ws.Range("E31").Value = "  -4.06%  "

# Row 32: update D, E
$ws.Range("D32").Value = "9.02"
$ws.Range("E32").Value = "  -4.01%  "

# Row 33: update E
$ws.Range("E33").Value = "  -0.02%  "

# Row 34: update D, E
$ws.Range("D34").Value = "4.92"
$ws.Range("E34").Value = "  -4.27%  "

# Row 35: update D, E
$ws.Range("D35").Value = "4.47"
$ws.Range("E35").Value = "  -3.64%  "

# Row 37: update D, E
$ws.Range("D37").Value = "0.0689"
$ws.Range("E37").Value = "  -5.57%  "

# Row 38: update D, E
$ws.Range("D38").Value = "16.17"
$ws.Range("E38").Value = "  -8.05%  "

# Row 39: update E
$ws.Range("E39").Value = "  -5.39%  "

# Row 40: update D, E
$ws.Range("D40").Value = "0.0987"
$ws.Range("E40").Value = "  -4.07%  "

# Row 41: update E
$ws.Range("E41").Value = "  -3.80%  "

# Row 42: update E
$ws.Range("E42").Value = "  -8.32%  "

# Row 43: update D, E
$ws.Range("D43").Value = "2.48"
$ws.Range("E43").Value = "  +0.89%  "

# Row 44: update D, E
$ws.Range("D44").Value = "1.959.81"
$ws.Range("E44").Value = "  -2.32%  "

# Row 45: update E
$ws.Range("E45").Value = "  -3.17%  "

# Row 46: update D, E
$ws.Range("D46").Value = "17.27"
$ws.Range("E46").Value = "  -8.10%  "

# Row 47: update D, E
$ws.Range("D47").Value = "9.54"
$ws.Range("E47").Value = "  -6.53%  "

# Row 48: update E
$ws.Range("E48").Value = "  -5.99%  "

# Row 49: update B, C, D, E
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.489.39"
$ws.Range("E49").Value = "  -2.57%  "

# Row 50: update B, C, D, E
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.79"
$ws.Range("E50").Value = "  -3.70%  "

# Row 51: update D, E
$ws.Range("D51").Value = "51.76"
$ws.Range("E51").Value = "  -7.30%  "

